$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.723.57"
$ws.Range("E2").Value = "  +2.37%  "
$ws.Range("D3").Value = "3.031.72"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "511.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.10%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.442"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.27%  "
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.111"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.368"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.64%  "
$ws.Range("D12").Value = "3.549.96"
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.82%  "
$ws.Range("E15").Value = "  +10.46%  "
$ws.Range("D16").Value = "57.718.49"
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.04%  "
$ws.Range("D18").Value = "3.026.96"
$ws.Range("E18").Value = "  +2.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "334.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.96%  "
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.500"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.169"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").Value = "0.0₃0929"
$ws.Range("E28").Value = "  +5.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.09%  "
$ws.Range("E31").Value = "  +4.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.73"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0687"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.61%  "
$ws.Range("D40").Value = "3.066.13"
$ws.Range("E40").Value = "  +2.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "2.315.64"
$ws.Range("E44").Value = "  +3.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.656"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.992"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.17%  "
